# Update "想去人数" (F column) values for several event rows across the
# "展览" (exhibitions), "演出" (performances) and "全部类型" (all types)
# sheets, per the source diff.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F6").Value = 130
$wsExhibition.Range("F9").Value = 328
$wsExhibition.Range("F13").Value = 11566

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 102

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F5").Value = 102
$wsAllTypes.Range("F8").Value = 130
$wsAllTypes.Range("F11").Value = 328
$wsAllTypes.Range("F15").Value = 11566
